# Adds the "ODI Bowling Extra" worksheet (MATCH_CODE / MAIDEN_OVERS /
# PERCENT_WICKETS_OF_ALL) as the 5th sheet, mirroring the structure of the
# existing "ODI Batting Extra" sheet.

$wb = $excel.ActiveWorkbook

# --- create the new sheet, positioned after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- copy the header formatting (bold / border / centered) from the
#     existing "ODI Batting Extra" header row so we reuse the same style ---
$template = $wb.Worksheets.Item("ODI Batting Extra")
$template.Range("A1:C1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats

# --- header row ---
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# --- data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL ---
# Force text formatting on the data range so values like "0" / "30.00%"
# are stored as text, matching the source data (not auto-converted to
# numbers/percentages).
$data = @(
    @("4232", "0", ""),
    @("4233", "", ""),
    @("4261", "0", ""),
    @("4264", "0", "30.00%"),
    @("4269", "0", ""),
    @("4271", "", ""),
    @("4272", "0", "10.00%"),
    @("4302", "", ""),
    @("4305", "0", ""),
    @("4309", "0", "10.00%"),
    @("4322", "0", ""),
    @("4331", "0", ""),
    @("4339", "1", ""),
    @("4350", "0", ""),
    @("4356", "", ""),
    @("4413", "1", "10.00%"),
    @("4414", "", ""),
    @("4417", "0", ""),
    @("4450", "", ""),
    @("4451", "0", "10.00%")
)

$row = 2
foreach ($rec in $data) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $rec[0]

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $rec[1]

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $rec[2]

    $row = $row + 1
}
